# Apply updated TPM-derived NATMI metrics (ligand/receptor expression,
# specificity, and edge-weight columns) to sheet1 of the workbook.
# Values below are the recomputed outputs after the TPM normalization
# script update described in the commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 51.282378
$ws.Range("H2").Value = 153.847134
$ws.Range("I2").Value = 0.6626772651633268
$ws.Range("J2").Value = 0.6626772651633267
$ws.Range("M2").Value = 361.2779286666667
$ws.Range("N2").Value = 1083.833786
$ws.Range("O2").Value = 0.9679392703861037
$ws.Range("P2").Value = 0.9679392703861038
$ws.Range("Q2").Value = 18527.19130094104
$ws.Range("R2").Value = 166744.7217084694
$ws.Range("S2").Value = 0.6414313485436491
$ws.Range("T2").Value = 0.6414313485436491

# Row 3
$ws.Range("G3").Value = 51.282378
$ws.Range("H3").Value = 153.847134
$ws.Range("I3").Value = 0.6626772651633268
$ws.Range("J3").Value = 0.6626772651633267
$ws.Range("O3").Value = 0.015995373883918
$ws.Range("P3").Value = 0.015995373883918
$ws.Range("Q3").Value = 306.165232617552
$ws.Range("R3").Value = 2755.487093557968
$ws.Range("S3").Value = 0.01059977062065968
$ws.Range("T3").Value = 0.01059977062065968

# Row 4
$ws.Range("G4").Value = 51.282378
$ws.Range("H4").Value = 153.847134
$ws.Range("I4").Value = 0.6626772651633268
$ws.Range("J4").Value = 0.6626772651633267
$ws.Range("M4").Value = 3.890485666666667
$ws.Range("N4").Value = 11.671457
$ws.Range("O4").Value = 0.01042342628440887
$ws.Range("P4").Value = 0.01042342628440887
$ws.Range("Q4").Value = 199.513356561582
$ws.Range("R4").Value = 1795.620209054238
$ws.Range("S4").Value = 0.006907367623783609
$ws.Range("T4").Value = 0.006907367623783609

# Row 5
$ws.Range("G5").Value = 51.282378
$ws.Range("H5").Value = 153.847134
$ws.Range("I5").Value = 0.6626772651633268
$ws.Range("J5").Value = 0.6626772651633267
$ws.Range("M5").Value = 2.105818666666667
$ws.Range("N5").Value = 6.317456
$ws.Range("O5").Value = 0.005641929445569353
$ws.Range("P5").Value = 0.005641929445569354
$ws.Range("Q5").Value = 107.991388863456
$ws.Range("R5").Value = 971.9224997711041
$ws.Range("S5").Value = 0.003738778375234343
$ws.Range("T5").Value = 0.003738778375234343

# Row 6
$ws.Range("I6").Value = 0.05445697206111803
$ws.Range("J6").Value = 0.05445697206111802
$ws.Range("M6").Value = 361.2779286666667
$ws.Range("N6").Value = 1083.833786
$ws.Range("O6").Value = 0.9679392703861037
$ws.Range("P6").Value = 0.9679392703861038
$ws.Range("Q6").Value = 1522.512981937999
$ws.Range("R6").Value = 13702.616837442
$ws.Range("S6").Value = 0.05271104180427501
$ws.Range("T6").Value = 0.05271104180427501

# Row 7
$ws.Range("I7").Value = 0.05445697206111803
$ws.Range("J7").Value = 0.05445697206111802
$ws.Range("O7").Value = 0.015995373883918
$ws.Range("P7").Value = 0.015995373883918
$ws.Range("S7").Value = 0.0008710596287036593
$ws.Range("T7").Value = 0.0008710596287036592

# Row 8
$ws.Range("I8").Value = 0.05445697206111803
$ws.Range("J8").Value = 0.05445697206111802
$ws.Range("M8").Value = 3.890485666666667
$ws.Range("N8").Value = 11.671457
$ws.Range("O8").Value = 0.01042342628440887
$ws.Range("P8").Value = 0.01042342628440887
$ws.Range("Q8").Value = 16.39545198735033
$ws.Range("R8").Value = 147.559067886153
$ws.Range("S8").Value = 0.0005676282339511773
$ws.Range("T8").Value = 0.0005676282339511773

# Row 9
$ws.Range("I9").Value = 0.05445697206111803
$ws.Range("J9").Value = 0.05445697206111802
$ws.Range("M9").Value = 2.105818666666667
$ws.Range("N9").Value = 6.317456
$ws.Range("O9").Value = 0.005641929445569353
$ws.Range("P9").Value = 0.005641929445569354
$ws.Range("Q9").Value = 8.874431575269332
$ws.Range("R9").Value = 79.869884177424
$ws.Range("S9").Value = 0.0003072423941881694
$ws.Range("T9").Value = 0.0003072423941881694

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2493196666666667
$ws.Range("H10").Value = 0.747959
$ws.Range("I10").Value = 0.003221739734029084
$ws.Range("J10").Value = 0.003221739734029083
$ws.Range("M10").Value = 361.2779286666667
$ws.Range("N10").Value = 1083.833786
$ws.Range("O10").Value = 0.9679392703861037
$ws.Range("P10").Value = 0.9679392703861038
$ws.Range("Q10").Value = 90.07369274919712
$ws.Range("R10").Value = 810.6632347427742
$ws.Range("S10").Value = 0.003118448407530031
$ws.Range("T10").Value = 0.003118448407530031

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2493196666666667
$ws.Range("H11").Value = 0.747959
$ws.Range("I11").Value = 0.003221739734029084
$ws.Range("J11").Value = 0.003221739734029083
$ws.Range("O11").Value = 0.015995373883918
$ws.Range("P11").Value = 0.015995373883918
$ws.Range("Q11").Value = 1.488484284818667
$ws.Range("R11").Value = 13.396358563368
$ws.Range("S11").Value = 0.00005153293160246972
$ws.Range("T11").Value = 0.00005153293160246971

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2493196666666667
$ws.Range("H12").Value = 0.747959
$ws.Range("I12").Value = 0.003221739734029084
$ws.Range("J12").Value = 0.003221739734029083
$ws.Range("M12").Value = 3.890485666666667
$ws.Range("N12").Value = 11.671457
$ws.Range("O12").Value = 0.01042342628440887
$ws.Range("P12").Value = 0.01042342628440887
$ws.Range("Q12").Value = 0.9699745895847779
$ws.Range("R12").Value = 8.729771306263
$ws.Range("S12").Value = 0.0000335815666252032
$ws.Range("T12").Value = 0.0000335815666252032

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2493196666666667
$ws.Range("H13").Value = 0.747959
$ws.Range("I13").Value = 0.003221739734029084
$ws.Range("J13").Value = 0.003221739734029083
$ws.Range("M13").Value = 2.105818666666667
$ws.Range("N13").Value = 6.317456
$ws.Range("O13").Value = 0.005641929445569353
$ws.Range("P13").Value = 0.005641929445569354
$ws.Range("Q13").Value = 0.5250220080337777
$ws.Range("R13").Value = 4.725198072304
$ws.Range("S13").Value = 0.00001817682827137946
$ws.Range("T13").Value = 0.00001817682827137946

# Row 14
$ws.Range("G14").Value = 21.640716
$ws.Range("H14").Value = 64.92214800000001
$ws.Range("I14").Value = 0.2796440230415261
$ws.Range("J14").Value = 0.279644023041526
$ws.Range("M14").Value = 361.2779286666667
$ws.Range("N14").Value = 1083.833786
$ws.Range("O14").Value = 0.9679392703861037
$ws.Range("P14").Value = 0.9679392703861038
$ws.Range("Q14").Value = 7818.313051343593
$ws.Range("R14").Value = 70364.81746209234
$ws.Range("S14").Value = 0.2706784316306495
$ws.Range("T14").Value = 0.2706784316306495

# Row 15
$ws.Range("G15").Value = 21.640716
$ws.Range("H15").Value = 64.92214800000001
$ws.Range("I15").Value = 0.2796440230415261
$ws.Range("J15").Value = 0.279644023041526
$ws.Range("O15").Value = 0.015995373883918
$ws.Range("P15").Value = 0.015995373883918
$ws.Range("Q15").Value = 129.199056411744
$ws.Range("R15").Value = 1162.791507705696
$ws.Range("S15").Value = 0.004473010702952189
$ws.Range("T15").Value = 0.004473010702952188

# Row 16
$ws.Range("G16").Value = 21.640716
$ws.Range("H16").Value = 64.92214800000001
$ws.Range("I16").Value = 0.2796440230415261
$ws.Range("J16").Value = 0.279644023041526
$ws.Range("M16").Value = 3.890485666666667
$ws.Range("N16").Value = 11.671457
$ws.Range("O16").Value = 0.01042342628440887
$ws.Range("P16").Value = 0.01042342628440887
$ws.Range("Q16").Value = 84.19289541440401
$ws.Range("R16").Value = 757.7360587296361
$ws.Range("S16").Value = 0.002914848860048884
$ws.Range("T16").Value = 0.002914848860048883

# Row 17
$ws.Range("G17").Value = 21.640716
$ws.Range("H17").Value = 64.92214800000001
$ws.Range("I17").Value = 0.2796440230415261
$ws.Range("J17").Value = 0.279644023041526
$ws.Range("M17").Value = 2.105818666666667
$ws.Range("N17").Value = 6.317456
$ws.Range("O17").Value = 0.005641929445569353
$ws.Range("P17").Value = 0.005641929445569354
$ws.Range("Q17").Value = 45.571423712832
$ws.Range("R17").Value = 410.142813415488
$ws.Range("S17").Value = 0.001577731847875461
$ws.Range("T17").Value = 0.00157773184787546
